$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.296772003173828
$ws.Range("B1").Value = 3.930430173873901
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.446487188339233
$ws.Range("E1").Value = 1.940884113311768
